$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data as per the Oct 25 2023 refresh.
# Force text number format before assignment so numeric-looking strings
# (e.g. "0.999", "1.20") remain stored as text, matching the source data,
# then restore the General format so styling is unaffected.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.677.73"
$ws.Range("D2").NumberFormat = "General"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("E2").NumberFormat = "General"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.797.19"
$ws.Range("D3").NumberFormat = "General"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E3").NumberFormat = "General"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").NumberFormat = "General"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("E4").NumberFormat = "General"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.71"
$ws.Range("D5").NumberFormat = "General"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("E5").NumberFormat = "General"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.555"
$ws.Range("D6").NumberFormat = "General"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.87%  "
$ws.Range("E6").NumberFormat = "General"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").NumberFormat = "General"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("E7").NumberFormat = "General"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.79"
$ws.Range("D8").NumberFormat = "General"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.56%  "
$ws.Range("E8").NumberFormat = "General"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("E9").NumberFormat = "General"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0676"
$ws.Range("D10").NumberFormat = "General"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("E10").NumberFormat = "General"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("D11").NumberFormat = "General"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("E11").NumberFormat = "General"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.056.66"
$ws.Range("D12").NumberFormat = "General"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("E12").NumberFormat = "General"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +10.07%  "
$ws.Range("E13").NumberFormat = "General"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.799.27"
$ws.Range("D14").NumberFormat = "General"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("E14").NumberFormat = "General"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.640"
$ws.Range("D15").NumberFormat = "General"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("E15").NumberFormat = "General"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.611.97"
$ws.Range("D16").NumberFormat = "General"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("E16").NumberFormat = "General"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("E17").NumberFormat = "General"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.71"
$ws.Range("D18").NumberFormat = "General"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("E18").NumberFormat = "General"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "256.32"
$ws.Range("D19").NumberFormat = "General"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("E19").NumberFormat = "General"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0779"
$ws.Range("D20").NumberFormat = "General"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.34%  "
$ws.Range("E20").NumberFormat = "General"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").NumberFormat = "General"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E21").NumberFormat = "General"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.47"
$ws.Range("D22").NumberFormat = "General"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.38%  "
$ws.Range("E22").NumberFormat = "General"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("E23").NumberFormat = "General"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.58%  "
$ws.Range("E24").NumberFormat = "General"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.35"
$ws.Range("D25").NumberFormat = "General"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("E25").NumberFormat = "General"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.53"
$ws.Range("D26").NumberFormat = "General"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("E26").NumberFormat = "General"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.15"
$ws.Range("D27").NumberFormat = "General"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.51%  "
$ws.Range("E27").NumberFormat = "General"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.47%  "
$ws.Range("E28").NumberFormat = "General"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").NumberFormat = "General"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("E29").NumberFormat = "General"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.81"
$ws.Range("D30").NumberFormat = "General"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.98%  "
$ws.Range("E30").NumberFormat = "General"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("E31").NumberFormat = "General"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.20"
$ws.Range("D32").NumberFormat = "General"

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("E32").NumberFormat = "General"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.61"
$ws.Range("D33").NumberFormat = "General"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("E33").NumberFormat = "General"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.90"
$ws.Range("D34").NumberFormat = "General"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.30%  "
$ws.Range("E34").NumberFormat = "General"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.455.64"
$ws.Range("D35").NumberFormat = "General"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.11%  "
$ws.Range("E35").NumberFormat = "General"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.78%  "
$ws.Range("E36").NumberFormat = "General"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.87%  "
$ws.Range("E37").NumberFormat = "General"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("E38").NumberFormat = "General"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "83.18"
$ws.Range("D39").NumberFormat = "General"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("E39").NumberFormat = "General"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.56%  "
$ws.Range("E40").NumberFormat = "General"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.33"
$ws.Range("D41").NumberFormat = "General"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.16%  "
$ws.Range("E41").NumberFormat = "General"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.905"
$ws.Range("D42").NumberFormat = "General"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("E42").NumberFormat = "General"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.09"
$ws.Range("D43").NumberFormat = "General"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("E43").NumberFormat = "General"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0508"
$ws.Range("D44").NumberFormat = "General"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.36%  "
$ws.Range("E44").NumberFormat = "General"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("E45").NumberFormat = "General"

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("B46").NumberFormat = "General"

$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C46").NumberFormat = "General"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.05"
$ws.Range("D46").NumberFormat = "General"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.14%  "
$ws.Range("E46").NumberFormat = "General"

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("B47").NumberFormat = "General"

$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("C47").NumberFormat = "General"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.955.74"
$ws.Range("D47").NumberFormat = "General"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("E47").NumberFormat = "General"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.03"
$ws.Range("D48").NumberFormat = "General"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("E48").NumberFormat = "General"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("D49").NumberFormat = "General"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("E49").NumberFormat = "General"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.56"
$ws.Range("D50").NumberFormat = "General"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.78%  "
$ws.Range("E50").NumberFormat = "General"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.41"
$ws.Range("D51").NumberFormat = "General"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.30%  "
$ws.Range("E51").NumberFormat = "General"
